$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for columns D (Fecha), J (Volumen), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# Each row below: RowIndex, D, J, L, M, P
$rowsData = @(
    ,@(2, 44307, 160, 10000, 10000, 3333)
    ,@(3, 44377, 16, 10500, 10250, 3417)
    ,@(4, 44356, 16, 10000, 10000, 3333)
    ,@(5, 44349, 12, 10000, 10000, 3333)
    ,@(6, 44266, 160, 10000, 10000, 3333)
    ,@(7, 44405, 16, 10500, 10250, 3417)
    ,@(8, 44363, 16, 10000, 10000, 3333)
    ,@(9, 44181, 10, 12000, 11000, 3667)
    ,@(10, 44328, 16, 10000, 10000, 3333)
    ,@(11, 44293, 16, 10000, 10000, 3333)
    ,@(12, 44195, 30, 10000, 10000, 3333)
    ,@(13, 44419, 16, 10000, 10000, 3333)
    ,@(14, 44370, 16, 10500, 10250, 3417)
    ,@(15, 44300, 16, 10000, 10000, 3333)
    ,@(16, 44384, 25, 10500, 10260, 3420)
    ,@(17, 44272, 70, 10000, 10000, 3333)
    ,@(18, 44433, 16, 10500, 10250, 3417)
    ,@(19, 44321, 25, 10000, 10000, 3333)
    ,@(20, 44335, 16, 10000, 10000, 3333)
    ,@(21, 44426, 16, 10500, 10250, 3417)
    ,@(22, 44398, 16, 10500, 10250, 3417)
    ,@(23, 44342, 17, 10000, 10000, 3333)
    ,@(24, 44391, 16, 10000, 10000, 3333)
    ,@(25, 44435, 16, 10500, 10250, 3417)
    ,@(26, 44279, 16, 10000, 10000, 3333)
    ,@(27, 44412, 25, 10500, 10260, 3420)
    ,@(28, 44314, 16, 10000, 10000, 3333)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J - Volumen
    $ws.Cells.Item($r, 12).Value = $row[3]  # L - Precio máximo
    $ws.Cells.Item($r, 13).Value = $row[4]  # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[5]  # P - Precio $/Kg
}
